$wb = $excel.ActiveWorkbook

# Both the Weekly Report and Monthly Report sheets receive identical edits:
#   - a handful of Quantity corrections
#   - updated contact info for the ingredient supplier ("Fresh Foods Inc.")
#   - a brand-new product row (PRD016 / Nestea (330ml))
$sheetNames = @("Weekly Report", "Monthly Report")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Quantity (column D) corrections ---
    $ws.Range("D2").Value = 39
    $ws.Range("D3").Value = 59
    $ws.Range("D4").Value = 70
    $ws.Range("D11").Value = 1
    $ws.Range("D16").Value = 24

    # --- Supplier contact-info corrections for Fresh Foods Inc. rows (6-10) ---
    # Phone numbers are all-digit strings, so they must be forced to Text
    # (otherwise Excel auto-coerces them into a number and drops the leading 0).
    $ws.Range("M6:M10").NumberFormat = "@"
    $ws.Range("M6:M10").Value = "09456123481"
    $ws.Range("M6:M10").ClearFormats()

    $ws.Range("N6:N10").Value = "info@tiamzonfoodco.com"
    $ws.Range("O6:O10").Value = "273 Harvard Avenue, Pasig"

    # --- New row 17: PRD016 / Nestea (330ml) ---
    $ws.Range("A17").Value = "PRD016"

    $ws.Range("B17").Value = 45465
    $ws.Range("B17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("C17").Value = "Nestea (330ml)"
    $ws.Range("D17").Value = 0
    $ws.Range("E17").Value = 10

    # Expiry_Date column stores plain text dates, not real date values.
    $ws.Range("F17").NumberFormat = "@"
    $ws.Range("F17").Value = "2025-03-01"
    $ws.Range("F17").ClearFormats()

    $ws.Range("G17").Value = "Low Stock"
    $ws.Range("H17").Value = "Beverage"
    $ws.Range("I17").Value = "Active"
    $ws.Range("J17").Value = 65
    $ws.Range("K17").Value = 40
    $ws.Range("L17").Value = "Global Food Distributors"

    $ws.Range("M17").NumberFormat = "@"
    $ws.Range("M17").Value = "09356789101"
    $ws.Range("M17").ClearFormats()

    $ws.Range("N17").Value = "info@globalfooddist.com"
    $ws.Range("O17").Value = "345 P. Burgos Street, Makati"
    $ws.Range("P17").Value = "Active"
}

Write-Output "edits applied"
